# Insert a new weekly price record for "Puerro" (Vega Modelo de Temuco)
# as row 74, pushing all the existing records (old rows 74-116) down by
# one row (new rows 75-117). Mirrors the author's commit: a new weekly
# observation was prepended to the time series kept in this sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Push existing data down to make room for the new record.
$ws.Rows("74:74").Insert()

# Populate the newly inserted row with the new observation.
$ws.Cells.Item(74, 1).Value  = 10
$ws.Cells.Item(74, 2).Value  = "Vega Modelo de Temuco"
$ws.Cells.Item(74, 3).Value  = "La Araucanía"
$ws.Cells.Item(74, 4).Value  = 44438
$ws.Cells.Item(74, 5).Value  = 9
$ws.Cells.Item(74, 6).Value  = 100112005
$ws.Cells.Item(74, 7).Value  = "Puerro"
$ws.Cells.Item(74, 8).Value  = "Azul de Maquehue"
$ws.Cells.Item(74, 9).Value  = "Primera"
$ws.Cells.Item(74, 10).Value = 50
$ws.Cells.Item(74, 11).Value = 8000
$ws.Cells.Item(74, 12).Value = 8000
$ws.Cells.Item(74, 13).Value = 8000
$ws.Cells.Item(74, 14).Value = "`$/docena de paquetes"
$ws.Cells.Item(74, 15).Value = "Provincia de Cautín"
$ws.Cells.Item(74, 16).Value = 667
$ws.Cells.Item(74, 17).Value = 12
$ws.Cells.Item(74, 18).Value = "Hortaliza"

# Give the new date cell the same date style ("s=2" numFmt) used by the
# rest of column D.
$ws.Cells.Item(74, 4).NumberFormat = $ws.Cells.Item(75, 4).NumberFormat
